$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.702.66'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '2.187.56'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'" + '292.86'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').Value = "'" + '86.22'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = "'" + '29.75'
$ws.Range('E10').Value = '  -3.76%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'" + '0.0773'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').Value = "'" + '6.37'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '2.530.21'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').Value = '2.233.64'
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = "'" + '13.57'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').Value = '39.604.05'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').Value = "'" + '11.12'
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('D22').Value = "'" + '5.70'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('D23').Value = "'" + '64.81'
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('D24').Value = "'" + '235.36'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('D27').Value = "'" + '1.79'
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('D28').Value = "'" + '22.28'
$ws.Range('E28').Value = '  -2.68%  '
$ws.Range('E29').Value = '  -3.65%  '
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').Value = "'" + '155.97'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').Value = "'" + '31.03'
$ws.Range('E32').Value = '  -6.73%  '
$ws.Range('D33').Value = "'" + '0.999'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = "'" + '2.77'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').Value = "'" + '0.0967'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('E40').Value = '  -7.60%  '
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('D42').Value = '2.105.67'
$ws.Range('E42').Value = '  +3.19%  '
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').Value = "'" + '17.14'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('D47').Value = "'" + '9.53'
$ws.Range('E47').Value = '  -5.85%  '
$ws.Range('D48').Value = "'" + '2.61'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').Value = '2.406.58'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('E50').Value = '  +3.52%  '
$ws.Range('D51').Value = "'" + '1.09'
$ws.Range('E51').Value = '  +0.62%  '
